$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 187, shifting the existing rows 187:286 down to 188:287.
$ws.Rows("187:187").Insert()

# Populate the newly inserted row 187 with the new record.
$ws.Range("A187").Value = 10
$ws.Range("B187").Value = 'Vega Modelo de Temuco'
$ws.Range("C187").Value = 'La Araucanía'
$ws.Range("D187").Value = 44193
$ws.Range("E187").Value = 9
$ws.Range("F187").Value = 'Fruta'
$ws.Range("G187").Value = 100103
$ws.Range("H187").Value = 'Frutos de hueso (carozo)'
$ws.Range("I187").Value = 100103006
$ws.Range("J187").Value = 'Nectarín'
$ws.Range("K187").Value = 'Super Queen'
$ws.Range("L187").Value = 'Primera'
$ws.Range("M187").Value = 210
$ws.Range("N187").Value = 18000
$ws.Range("O187").Value = 18000
$ws.Range("P187").Value = 18000
$ws.Range("Q187").Value = '$/bandeja 10 kilos granel'
$ws.Range("R187").Value = "Región de O'Higgins"
$ws.Range("S187").Value = 1800
$ws.Range("T187").Value = 10
